# Daily attendance processing - 2025-11-29 19:21:37
# Reverse the order of the comma-separated "Recorded By" entries in column G
# (e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count + $ws.UsedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = [string]$val -split ",\s*"
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $cell.Value2 = $reversed -join ", "
        }
    }
}
